$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1802.7
$ws.Range("I4").Value = 904.2857
$ws.Range("K4").Value = 904.2857
$ws.Range("M4").Value = -790.2857
$ws.Range("H9").Value = 459.8
$ws.Range("I9").Value = 341.58823
$ws.Range("J9").Value = 1129.6666
$ws.Range("K9").Value = 341.58823
$ws.Range("L9").Value = 1129.6666
$ws.Range("M9").Value = -172.58823
$ws.Range("N9").Value = -1467.6666
$ws.Range("H40").Value = 5467.758
$ws.Range("I40").Value = 4195.25
$ws.Range("J40").Value = 8861.111000000001
$ws.Range("K40").Value = 4195.25
$ws.Range("L40").Value = 8861.111000000001
$ws.Range("M40").Value = -4020.25
$ws.Range("N40").Value = -9211.111000000001
$ws.Range("H43").Value = 6925
$ws.Range("J43").Value = 4800
$ws.Range("L43").Value = 4800
$ws.Range("N43").Value = -4938
$ws.Range("H80").Value = 912.1818
$ws.Range("I80").Value = 366.2
$ws.Range("J80").Value = 1367.1666
$ws.Range("K80").Value = 1098.6
$ws.Range("L80").Value = 4101.4998
$ws.Range("M80").Value = -100.5999999999999
$ws.Range("N80").Value = -6097.4998
$ws.Range("H83").Value = 912.1818
$ws.Range("I83").Value = 366.2
$ws.Range("J83").Value = 1367.1666
$ws.Range("K83").Value = 3295.8
$ws.Range("L83").Value = 12304.4994
$ws.Range("M83").Value = 1696.2
$ws.Range("N83").Value = -22288.4994
$ws.Range("H94").Value = 498.33334
$ws.Range("I94").Value = 498.33334
$ws.Range("K94").Value = 498.33334
$ws.Range("M94").Value = -47.33334000000002
$ws.Range("H97").Value = 3881.6667
$ws.Range("J97").Value = 4419.2856
$ws.Range("L97").Value = 13257.8568
$ws.Range("N97").Value = -14249.8568
$ws.Range("H99").Value = 3534.0833
$ws.Range("I99").Value = 2929.5
$ws.Range("J99").Value = 4743.25
$ws.Range("K99").Value = 8788.5
$ws.Range("L99").Value = 14229.75
$ws.Range("M99").Value = -7290.5
$ws.Range("N99").Value = -17225.75
$ws.Range("H107").Value = 1109.0454
$ws.Range("I107").Value = 1038
$ws.Range("J107").Value = 1559
$ws.Range("K107").Value = 1038
$ws.Range("L107").Value = 1559
$ws.Range("M107").Value = 882
$ws.Range("N107").Value = -5399
$ws.Range("H138").Value = 2491.46
$ws.Range("I138").Value = 1800.3043
$ws.Range("J138").Value = 3080.2222
$ws.Range("K138").Value = 5400.9129
$ws.Range("L138").Value = 9240.6666
$ws.Range("M138").Value = -260.9129000000003
$ws.Range("N138").Value = -19520.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -887
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 415
$ws.Range("I5").Value = 453.33334
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 453.33334
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -341.33334
$ws.Range("N5").Value = -524
$ws.Range("H116").Value = 1000
$ws.Range("I116").Value = 1000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1294
$ws.Range("N116").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -886
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 415
$ws.Range("I4").Value = 453.33334
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 453.33334
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -338.33334
$ws.Range("N4").Value = -530
$ws.Range("H22").Value = 1099.9286
$ws.Range("I22").Value = 949.875
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 949.875
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = -776.875
$ws.Range("N22").Value = -1646
$ws.Range("H86").Value = 1104.3928
$ws.Range("J86").Value = 2503.5
$ws.Range("L86").Value = 2503.5
$ws.Range("N86").Value = -4749.5
$ws.Range("H89").Value = 1104.3928
$ws.Range("J89").Value = 2503.5
$ws.Range("L89").Value = 12517.5
$ws.Range("N89").Value = -23749.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 691.55554
$ws.Range("I7").Value = 787.58826
$ws.Range("J7").Value = 528.3
$ws.Range("K7").Value = 787.58826
$ws.Range("L7").Value = 528.3
$ws.Range("M7").Value = -674.58826
$ws.Range("N7").Value = -754.3
$ws.Range("H16").Value = 698.1429000000001
$ws.Range("I16").Value = 698.1429000000001
$ws.Range("K16").Value = 698.1429000000001
$ws.Range("M16").Value = -411.1429000000001
$ws.Range("H31").Value = 8201.862999999999
$ws.Range("I31").Value = 4527.7334
$ws.Range("J31").Value = 10102.275
$ws.Range("K31").Value = 4527.7334
$ws.Range("L31").Value = 10102.275
$ws.Range("M31").Value = -4232.7334
$ws.Range("N31").Value = -10692.275
$ws.Range("H34").Value = 8201.862999999999
$ws.Range("I34").Value = 4527.7334
$ws.Range("J34").Value = 10102.275
$ws.Range("K34").Value = 4527.7334
$ws.Range("L34").Value = 10102.275
$ws.Range("M34").Value = -4325.7334
$ws.Range("N34").Value = -10506.275
$ws.Range("H58").Value = 3831.077
$ws.Range("I58").Value = 2780.7
$ws.Range("J58").Value = 7332.3335
$ws.Range("K58").Value = 2780.7
$ws.Range("L58").Value = 7332.3335
$ws.Range("M58").Value = -2577.7
$ws.Range("N58").Value = -7738.3335
$ws.Range("H59").Value = 105749
$ws.Range("J59").Value = 105749
$ws.Range("L59").Value = 105749
$ws.Range("N59").Value = -108039
$ws.Range("H113").Value = 698.1429000000001
$ws.Range("I113").Value = 698.1429000000001
$ws.Range("K113").Value = 698.1429000000001
$ws.Range("M113").Value = 1471.8571
$ws.Range("H136").Value = 3831.077
$ws.Range("I136").Value = 2780.7
$ws.Range("J136").Value = 7332.3335
$ws.Range("K136").Value = 8342.099999999999
$ws.Range("L136").Value = 21997.0005
$ws.Range("M136").Value = -5792.099999999999
$ws.Range("N136").Value = -27097.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 495.83334
$ws.Range("H117").Value = 2216.5
$ws.Range("J117").Value = 2324.75
$ws.Range("L117").Value = 6974.25
$ws.Range("N117").Value = -13858.25
$ws.Range("H131").Value = 4580.4443
$ws.Range("I131").Value = 2032.25
$ws.Range("J131").Value = 6619
$ws.Range("K131").Value = 6096.75
$ws.Range("L131").Value = 19857
$ws.Range("M131").Value = -1056.75
$ws.Range("N131").Value = -29937

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5020
$ws.Range("I80").Value = 5026.6665
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 5026.6665
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -4028.6665
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 5020
$ws.Range("I83").Value = 5026.6665
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 25133.3325
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -20141.3325
$ws.Range("N83").Value = -34984
$ws.Range("H96").Value = 16500
$ws.Range("J96").Value = 16500
$ws.Range("L96").Value = 16500
$ws.Range("N96").Value = -21992
$ws.Range("H101").Value = 63133.332
$ws.Range("J101").Value = 63133.332
$ws.Range("L101").Value = 63133.332
$ws.Range("N101").Value = -69623.33199999999
$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82080
$ws.Range("H132").Value = 7246
$ws.Range("I132").Value = 5215.4
$ws.Range("K132").Value = 15646.2
$ws.Range("M132").Value = -13116.2
$ws.Range("H141").Value = 56247.25
$ws.Range("J141").Value = 56247.25
$ws.Range("L141").Value = 56247.25
$ws.Range("N141").Value = -66607.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2016
$ws.Range("I16").Value = 2016
$ws.Range("K16").Value = 2016
$ws.Range("M16").Value = -1846
$ws.Range("H22").Value = 1661.409
$ws.Range("I22").Value = 1770.2307
$ws.Range("K22").Value = 1770.2307
$ws.Range("M22").Value = -1475.2307
$ws.Range("H27").Value = 1661.409
$ws.Range("I27").Value = 1770.2307
$ws.Range("K27").Value = 1770.2307
$ws.Range("M27").Value = -1663.2307
$ws.Range("H61").Value = 8122.727
$ws.Range("I61").Value = 3892.8
$ws.Range("J61").Value = 11647.667
$ws.Range("K61").Value = 3892.8
$ws.Range("L61").Value = 11647.667
$ws.Range("M61").Value = -3690.8
$ws.Range("N61").Value = -12051.667
$ws.Range("H82").Value = 2158.5715
$ws.Range("I82").Value = 1207.7142
$ws.Range("K82").Value = 1207.7142
$ws.Range("M82").Value = -846.7141999999999
$ws.Range("H85").Value = 2158.5715
$ws.Range("I85").Value = 1207.7142
$ws.Range("K85").Value = 1207.7142
$ws.Range("M85").Value = 40.28580000000011
$ws.Range("H113").Value = 8122.727
$ws.Range("I113").Value = 3892.8
$ws.Range("J113").Value = 11647.667
$ws.Range("K113").Value = 3892.8
$ws.Range("L113").Value = 11647.667
$ws.Range("M113").Value = -1722.8
$ws.Range("N113").Value = -15987.667
$ws.Range("H132").Value = 9452.441000000001
$ws.Range("I132").Value = 8086.0454
$ws.Range("J132").Value = 11957.5
$ws.Range("K132").Value = 24258.1362
$ws.Range("L132").Value = 35872.5
$ws.Range("M132").Value = -21728.1362
$ws.Range("N132").Value = -40932.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1182.9556
$ws.Range("I107").Value = 1257.6129
$ws.Range("J107").Value = 1017.6429
$ws.Range("K107").Value = 3772.8387
$ws.Range("L107").Value = 3052.9287
$ws.Range("M107").Value = -1852.8387
$ws.Range("N107").Value = -6892.9287
$ws.Range("H132").Value = 2340.7256
$ws.Range("I132").Value = 1970.9333
$ws.Range("K132").Value = 5912.7999
$ws.Range("M132").Value = -3382.7999
